$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "58.810.83"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.654.54"
$ws.Range("E3").Value = "  +3.65%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.04"
$ws.Range("E5").Value = "  +1.42%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.77"
$ws.Range("E6").Value = "  +0.43%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.34%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +2.22%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.682.50"
$ws.Range("E9").Value = "  +4.71%  "

# Row 10 - Toncoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.20"
$ws.Range("E10").Value = "  -0.11%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +4.97%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.55%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.92%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.121.94"
$ws.Range("E14").Value = "  +3.89%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "58.799.33"
$ws.Range("E15").Value = "  -0.10%  "

# Row 16 - Avalanche
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.95"
$ws.Range("E16").Value = "  +1.85%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +2.04%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.681.41"
$ws.Range("E18").Value = "  +4.78%  "

# Row 19 - BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.48"
$ws.Range("E19").Value = "  +4.72%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +0.34%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  +3.42%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +3.44%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.16%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.87"
$ws.Range("E24").Value = "  +2.56%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  +3.28%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "2.785.56"
$ws.Range("E26").Value = "  +4.42%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.49%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +1.64%  "

# Row 29 & 30 swapped: InternetComputer(DFINITY) <-> PEPE
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0807"
$ws.Range("E29").Value = "  +4.12%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.21"
$ws.Range("E30").Value = "  +5.10%  "

# Row 31 - USDe
$ws.Range("E31").Value = "  -0.29%  "

# Row 32 - Aptos
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.43"
$ws.Range("E32").Value = "  +10.72%  "

# Row 33 - EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.94"
$ws.Range("E33").Value = "  +2.05%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  +2.28%  "

# Row 35 - Monero
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.96"
$ws.Range("E35").Value = "  +0.21%  "

# Row 36 - SuiNetwork
$ws.Range("E36").Value = "  +14.35%  "

# Row 37 - NEARProtocol
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.02"
$ws.Range("E37").Value = "  +2.85%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  +3.25%  "

# Row 39 - OKB
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.71"
$ws.Range("E39").Value = "  +2.30%  "

# Row 40 - Fetch.AI
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.841"
$ws.Range("E40").Value = "  +2.00%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +5.50%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  +1.40%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  +1.98%  "

# Row 44 - Bittensor
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "278.39"
$ws.Range("E44").Value = "  -2.92%  "

# Row 45 - FirstDigitalUSD
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.992"
$ws.Range("E45").Value = "  -0.48%  "

# Row 46 - Stellar
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0979"
$ws.Range("E46").Value = "  +0.09%  "

# Row 47 - EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.73"
$ws.Range("E47").Value = "  +6.05%  "

# Row 48 - Hedera
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0532"
$ws.Range("E48").Value = "  +0.50%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +1.93%  "

# Row 50 & 51 swapped: RenderToken <-> Maker
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.002.55"
$ws.Range("E50").Value = "  +4.70%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.72"
$ws.Range("E51").Value = "  +4.36%  "
